$p = $ppt.ActivePresentation

# Locate the shape (normally "TextBox 2" on slide 1) holding the GitHub intro
# sentence, searching every slide defensively in case layout/order differs.
$sh = $null
for ($si = 1; $si -le $p.Slides.Count -and $sh -eq $null; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $cand = $slide.Shapes.Item($i)
        if ($cand.HasTextFrame -and $cand.TextFrame.HasText) {
            if ($cand.TextFrame.TextRange.Text -like "*는 버전 제어 및 공동 작업을 위한 코드 호스팅 플랫폼입니다*") {
                $sh = $cand
            }
        }
    }
}

$tr = $sh.TextFrame.TextRange
$full = $tr.Text

# The original sentence run reads:
#   "는 버전 제어 및 공동 작업을 위한 코드 호스팅 플랫폼입니다"
# It must be split into four runs:
#   "는 버전 제어 및 공동 작업을 위한 코드 "  (unchanged formatting)
#   "호스팅"
#   " "
#   "플랫폼 입니다"   (note the added space before "입니다")

$part1 = "는 버전 제어 및 공동 작업을 위한 코드 "
$part2 = "호스팅"
$part3 = " "
$oldPart4 = "플랫폼입니다"
$newPart4 = "플랫폼 입니다"

$startIdx = $full.IndexOf($part1)  # 0-based
$start1 = $startIdx + 1            # 1-based start of part1
$start2 = $start1 + $part1.Length
$start3 = $start2 + $part2.Length
$start4 = $start3 + $part3.Length

# Apply edits back-to-front so earlier offsets stay valid as the text length changes.
$c4 = $tr.Characters($start4, $oldPart4.Length)
$c4.Text = $newPart4

$c3 = $tr.Characters($start3, $part3.Length)
$c3.Text = $part3

$c2 = $tr.Characters($start2, $part2.Length)
$c2.Text = $part2
